# Swap the two theme palettes used by this deck (ppt/theme/theme1.xml and
# ppt/theme/theme2.xml): theme1.xml ("Integral" / Red Violet clrScheme)
# becomes the "Office Theme" palette, theme2.xml ("Office Theme") becomes
# the "Integral" (Red Violet) palette. Font scheme and format scheme were
# already identical between the two theme parts, so only the 12 scheme
# colors (and, best-effort, the theme/clrScheme display names) need to
# move.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

# Best-effort: rename the theme / color scheme to match the incoming
# "Office Theme" palette (PowerPoint shows these in the Design gallery).
try { $theme.Name = "Office Theme" } catch {}
try { $colors.Name = "Office" } catch {}

# Office Theme color scheme (target palette for theme1.xml), expressed as
# packed BGR integers (VBA/PowerPoint RGB long values = R + G*256 + B*65536)
# to match the .RGB property PowerPoint COM exposes on ThemeColor objects.
# Order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$colors.Item(1).RGB  = 0        # dk1      000000
$colors.Item(2).RGB  = 16777215 # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388  # dk2      44546A
$colors.Item(4).RGB  = 15132391 # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939 # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501  # accent2  ED7D31
$colors.Item(7).RGB  = 10855845 # accent3  A5A5A5
$colors.Item(8).RGB  = 49407    # accent4  FFC000
$colors.Item(9).RGB  = 12874308 # accent5  4472C4
$colors.Item(10).RGB = 4697456  # accent6  70AD47
$colors.Item(11).RGB = 12673797 # hlink    0563C1
$colors.Item(12).RGB = 7491477  # folHlink 954F72
